$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1808118081180812
$ws.Range("C2").Value = 0.5977859778597786
$ws.Range("P2").Value = 0.1180811808118081
$ws.Range("S2").Value = 0.1033210332103321
$ws.Range("B3").Value = 0.01754385964912281
$ws.Range("C3").Value = 0.04093567251461988
$ws.Range("P3").Value = 0.7660818713450293
$ws.Range("S3").Value = 0.1754385964912281
$ws.Range("J4").Value = 0.02173913043478261
$ws.Range("P4").Value = 0.7391304347826086
$ws.Range("S4").Value = 0.2391304347826087
$ws.Range("P5").Value = 0.6666666666666666
$ws.Range("S5").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.1055045871559633
$ws.Range("D6").Value = 0.02293577981651376
$ws.Range("F6").Value = 0.03669724770642202
$ws.Range("J6").Value = 0.2018348623853211
$ws.Range("O6").Value = 0.01376146788990826
$ws.Range("Q6").Value = 0.1605504587155963
$ws.Range("R6").Value = 0.1192660550458716
$ws.Range("S6").Value = 0.3394495412844037
$ws.Range("B7").Value = 0.1216216216216216
$ws.Range("D7").Value = 0.03378378378378379
$ws.Range("F7").Value = 0.07432432432432433
$ws.Range("J7").Value = 0.08783783783783784
$ws.Range("O7").Value = 0.01351351351351351
$ws.Range("Q7").Value = 0.1891891891891892
$ws.Range("R7").Value = 0.1013513513513514
$ws.Range("S7").Value = 0.3783783783783784
$ws.Range("B8").Value = 0.07947019867549669
$ws.Range("D8").Value = 0.02207505518763797
$ws.Range("E8").Value = 0.002207505518763797
$ws.Range("F8").Value = 0.0640176600441501
$ws.Range("J8").Value = 0.1258278145695364
$ws.Range("O8").Value = 0.02207505518763797
$ws.Range("Q8").Value = 0.17439293598234
$ws.Range("R8").Value = 0.09050772626931568
$ws.Range("S8").Value = 0.4194260485651214
$ws.Range("B9").Value = 0.08144796380090498
$ws.Range("D9").Value = 0.02262443438914027
$ws.Range("F9").Value = 0.05882352941176471
$ws.Range("J9").Value = 0.1493212669683258
$ws.Range("O9").Value = 0.01357466063348416
$ws.Range("Q9").Value = 0.1764705882352941
$ws.Range("R9").Value = 0.05429864253393665
$ws.Range("S9").Value = 0.4434389140271493
$ws.Range("B10").Value = 0.09581320450885668
$ws.Range("D10").Value = 0.01932367149758454
$ws.Range("E10").Value = 0.00322061191626409
$ws.Range("F10").Value = 0.08293075684380032
$ws.Range("J10").Value = 0.1207729468599034
$ws.Range("O10").Value = 0.01288244766505636
$ws.Range("Q10").Value = 0.2181964573268921
$ws.Range("R10").Value = 0.1014492753623188
$ws.Range("S10").Value = 0.3454106280193237
$ws.Range("G11").Value = 0.1391304347826087
$ws.Range("J11").Value = 0.09130434782608696
$ws.Range("K11").Value = 0.1782608695652174
$ws.Range("L11").Value = 0.5782608695652174
$ws.Range("S11").Value = 0.01304347826086956
$ws.Range("G12").Value = 0.7194244604316546
$ws.Range("J12").Value = 0.1942446043165468
$ws.Range("L12").Value = 0.05035971223021583
$ws.Range("S12").Value = 0.03597122302158273
$ws.Range("F15").Value = 0.00966183574879227
$ws.Range("H15").Value = 0.1594202898550725
$ws.Range("I15").Value = 0.07246376811594203
$ws.Range("J15").Value = 0.4009661835748792
$ws.Range("K15").Value = 0.07246376811594203
$ws.Range("N15").Value = 0.004830917874396135
$ws.Range("O15").Value = 0.06763285024154589
$ws.Range("S15").Value = 0.21256038647343
$ws.Range("F16").Value = 0.01587301587301587
$ws.Range("H16").Value = 0.1904761904761905
$ws.Range("I16").Value = 0.1111111111111111
$ws.Range("J16").Value = 0.3862433862433862
$ws.Range("K16").Value = 0.08465608465608465
$ws.Range("M16").Value = 0.01058201058201058
$ws.Range("O16").Value = 0.04232804232804233
$ws.Range("S16").Value = 0.1587301587301587
$ws.Range("F17").Value = 0.01118568232662192
$ws.Range("H17").Value = 0.1677852348993289
$ws.Range("I17").Value = 0.1208053691275168
$ws.Range("J17").Value = 0.4630872483221476
$ws.Range("K17").Value = 0.07158836689038031
$ws.Range("M17").Value = 0.01565995525727069
$ws.Range("N17").Value = 0.002237136465324385
$ws.Range("O17").Value = 0.04250559284116331
$ws.Range("S17").Value = 0.1051454138702461
$ws.Range("F18").Value = 0.0184331797235023
$ws.Range("H18").Value = 0.1751152073732719
$ws.Range("I18").Value = 0.09216589861751152
$ws.Range("J18").Value = 0.3870967741935484
$ws.Range("K18").Value = 0.08755760368663594
$ws.Range("M18").Value = 0.009216589861751152
$ws.Range("O18").Value = 0.06912442396313365
$ws.Range("S18").Value = 0.1612903225806452
$ws.Range("F19").Value = 0.01254180602006689
$ws.Range("H19").Value = 0.2282608695652174
$ws.Range("I19").Value = 0.09448160535117058
$ws.Range("J19").Value = 0.3787625418060201
$ws.Range("K19").Value = 0.08779264214046822
$ws.Range("M19").Value = 0.01923076923076923
$ws.Range("N19").Value = 0.0008361204013377926
$ws.Range("S19").Value = 0.1011705685618729
